$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the variable description in C12: the "Chronic Illness" concept's
# variable used to describe "chronic illnesses" generically; it now
# describes the Hepatitus B & C proxy measure used for that concept.
$ws.Range("C12").Value = "The presence of Hepatitus B & C within a household (proxy for chronic illnesses)."

# Reflect the author's updated selection in the sheet (C12 selected as the
# active cell, matching the saved view state).
$ws.Range("C12").Select()
